$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 249.5
$ws.Range("I12").Value = 299
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 299
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = -129
$ws.Range("N12").Value = -540
$ws.Range("H74").Value = 3250.5
$ws.Range("I74").Value = 3501.5
$ws.Range("J74").Value = 2999.5
$ws.Range("K74").Value = 3501.5
$ws.Range("L74").Value = 2999.5
$ws.Range("M74").Value = -2565.5
$ws.Range("N74").Value = -4871.5
$ws.Range("H77").Value = 3250.5
$ws.Range("I77").Value = 3501.5
$ws.Range("J77").Value = 2999.5
$ws.Range("K77").Value = 17507.5
$ws.Range("L77").Value = 14997.5
$ws.Range("M77").Value = -12827.5
$ws.Range("N77").Value = -24357.5
$ws.Range("H116").Value = 3640.6
$ws.Range("J116").Value = 5068
$ws.Range("L116").Value = 5068
$ws.Range("N116").Value = -11952
$ws.Range("H132").Value = 7253679.5
$ws.Range("I132").Value = 8776725
$ws.Range("K132").Value = 26330175
$ws.Range("M132").Value = -26327645
$ws.Range("H138").Value = 1827.67
$ws.Range("I138").Value = 1059.2858
$ws.Range("J138").Value = 2031.9241
$ws.Range("K138").Value = 3177.8574
$ws.Range("L138").Value = 6095.7723
$ws.Range("M138").Value = 1962.1426
$ws.Range("N138").Value = -16375.7723

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1083.6666
$ws.Range("I2").Value = 826.1875
$ws.Range("J2").Value = 1598.625
$ws.Range("K2").Value = 826.1875
$ws.Range("L2").Value = 1598.625
$ws.Range("M2").Value = -713.1875
$ws.Range("N2").Value = -1824.625
$ws.Range("H32").Value = 7524.828
$ws.Range("I32").Value = 5764.4824
$ws.Range("K32").Value = 5764.4824
$ws.Range("M32").Value = -5477.4824
$ws.Range("H63").Value = 32260870
$ws.Range("I63").Value = 2121.7727
$ws.Range("K63").Value = 2121.7727
$ws.Range("M63").Value = -1435.7727
$ws.Range("H66").Value = 32260870
$ws.Range("I66").Value = 2121.7727
$ws.Range("K66").Value = 10608.8635
$ws.Range("M66").Value = -7176.863499999999
$ws.Range("H74").Value = 3062.6667
$ws.Range("I74").Value = 2225.125
$ws.Range("J74").Value = 3732.7
$ws.Range("K74").Value = 2225.125
$ws.Range("L74").Value = 3732.7
$ws.Range("M74").Value = -1351.125
$ws.Range("N74").Value = -5480.7
$ws.Range("H77").Value = 3062.6667
$ws.Range("I77").Value = 2225.125
$ws.Range("J77").Value = 3732.7
$ws.Range("K77").Value = 11125.625
$ws.Range("L77").Value = 18663.5
$ws.Range("M77").Value = -6757.625
$ws.Range("N77").Value = -27399.5
$ws.Range("H116").Value = 1083.6666
$ws.Range("I116").Value = 826.1875
$ws.Range("J116").Value = 1598.625
$ws.Range("K116").Value = 826.1875
$ws.Range("L116").Value = 1598.625
$ws.Range("M116").Value = 1467.8125
$ws.Range("N116").Value = -6186.625
$ws.Range("H132").Value = 2283.2156
$ws.Range("I132").Value = 1540.5883
$ws.Range("J132").Value = 3768.4707
$ws.Range("K132").Value = 4621.7649
$ws.Range("L132").Value = 11305.4121
$ws.Range("M132").Value = -2091.7649
$ws.Range("N132").Value = -16365.4121

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1083.6666
$ws.Range("I3").Value = 826.1875
$ws.Range("J3").Value = 1598.625
$ws.Range("K3").Value = 826.1875
$ws.Range("L3").Value = 1598.625
$ws.Range("M3").Value = -712.1875
$ws.Range("N3").Value = -1826.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 12600500
$ws.Range("I6").Value = 12600500
$ws.Range("K6").Value = 12600500
$ws.Range("M6").Value = -12600387
$ws.Range("H107").Value = 747.7619
$ws.Range("I107").Value = 391.86667
$ws.Range("K107").Value = 391.86667
$ws.Range("M107").Value = 1528.13333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 127431
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 127431
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 382293
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -387161
$ws.Range("H107").Value = 3454.9697
$ws.Range("I107").Value = 532.25
$ws.Range("J107").Value = 4390.24
$ws.Range("K107").Value = 1596.75
$ws.Range("L107").Value = 13170.72
$ws.Range("M107").Value = 323.25
$ws.Range("N107").Value = -17010.72
$ws.Range("H113").Value = 647.72095
$ws.Range("I113").Value = 568.7692
$ws.Range("K113").Value = 1706.3076
$ws.Range("M113").Value = 463.6924000000001
$ws.Range("H123").Value = 2980.3635
$ws.Range("I123").Value = 2972.5
$ws.Range("J123").Value = 2984.8572
$ws.Range("K123").Value = 8917.5
$ws.Range("L123").Value = 8954.571599999999
$ws.Range("M123").Value = -6467.5
$ws.Range("N123").Value = -13854.5716
$ws.Range("H131").Value = 25644482
$ws.Range("I131").Value = 100000400
$ws.Range("J131").Value = 4510.1377
$ws.Range("K131").Value = 300001200
$ws.Range("L131").Value = 13530.4131
$ws.Range("M131").Value = -299996160
$ws.Range("N131").Value = -23610.4131

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1431.4667
$ws.Range("I102").Value = 1472.8462
$ws.Range("K102").Value = 1472.8462
$ws.Range("M102").Value = 149.1538
$ws.Range("H109").Value = 6999.7144
$ws.Range("J109").Value = 6999.7144
$ws.Range("L109").Value = 6999.7144
$ws.Range("N109").Value = -9079.714400000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1237.9286
$ws.Range("I61").Value = 1143
$ws.Range("J61").Value = 1364.5
$ws.Range("K61").Value = 1143
$ws.Range("L61").Value = 1364.5
$ws.Range("M61").Value = -941
$ws.Range("N61").Value = -1768.5
$ws.Range("H113").Value = 1237.9286
$ws.Range("I113").Value = 1143
$ws.Range("J113").Value = 1364.5
$ws.Range("K113").Value = 1143
$ws.Range("L113").Value = 1364.5
$ws.Range("M113").Value = 1027
$ws.Range("N113").Value = -5704.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 26500
$ws.Range("J46").Value = 26500
$ws.Range("L46").Value = 26500
$ws.Range("N46").Value = -26962
$ws.Range("H132").Value = 2124.9707
$ws.Range("I132").Value = 1983.8518
$ws.Range("J132").Value = 2669.2856
$ws.Range("K132").Value = 5951.555399999999
$ws.Range("L132").Value = 8007.8568
$ws.Range("M132").Value = -3421.555399999999
$ws.Range("N132").Value = -13067.8568
$ws.Range("H134").Value = 26500
$ws.Range("J134").Value = 26500
$ws.Range("L134").Value = 79500
$ws.Range("N134").Value = -84570
